# Apply the changes described by the diff:
#  1. Insert a new "Player Info" worksheet before the existing "ODI Batting" sheet.
#  2. Populate "Player Info" with player metadata (ID, NAME, BATTING_HAND, BOWL_STYLE).
#  3. On "ODI Batting": rename column D header MATCH_CARD_LINK -> MATCH_CODE and
#     replace the full scorecard URLs with just the numeric match code.

$wb = $excel.ActiveWorkbook

# --- Step 1: locate the existing "ODI Batting" sheet and insert a new sheet before it ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# NOTE: worksheet references returned/held by this runtime behave positionally,
# so after inserting a sheet (which shifts indices), re-resolve "ODI Batting" by
# name rather than reusing the handle obtained before the insertion.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- Step 2: fill in the "Player Info" sheet ---
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A2").Value = "'3759"
$playerInfo.Range("B2").Value = "Wriddhiman Prasanta Saha"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

$playerInfo.Range("A1").Select()

# --- Step 3: update the "ODI Batting" sheet's MATCH_CARD_LINK column ---
$battingSheet.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{
    2 = "3204"
    3 = "3206"
    4 = "3208"
    5 = "3643"
    6 = "3644"
    7 = "3645"
    8 = "3688"
    9 = "3689"
    10 = "3692"
}

foreach ($row in $matchCodes.Keys) {
    $battingSheet.Range("D$row").Value = "'" + $matchCodes[$row]
}
